# Generate Report for Handoff
# Rewrites the localization-status report: the source .md handoff file got
# a new GUID (44f387ea-1a36-439f-b416-767111d2633c), its handoff timestamps
# advance by ~21s, and two new dependent .png assets are handed off
# alongside it. Overview gains two new rows (one per .png) and the zh-cn /
# de-de detail sheets gain matching rows + hyperlinks.

$wb = $excel.ActiveWorkbook

$mdGuid   = "44f387ea-1a36-439f-b416-767111d2633c"
$mdFile   = "$mdGuid.md"
$xlfHash  = "9aa9b95d1ecd75a7f1c0601a8de955ddedbd6601"
$png1     = "5a99fcbf-8be2-4b16-8e5f-6fc128856d9a.png"
$png2     = "9a6b5d32-a4c1-461a-8210-0c15911e2e96.png"
$png1Hb   = "2529ea6593eb99251f354f54289599783d2e2d0a.png"
$png2Hb   = "c2b351324bfe3a92d9c60d61b9117a0d4d3c4d6c.png"

$zhXlf    = "$mdGuid.$xlfHash.zh-cn.xlf"
$deXlf    = "$mdGuid.$xlfHash.de-de.xlf"

$overviewDate = "2016-46-19 00:46:56"
$zhDate       = "2016-03-19 00:46:53"
$deDate       = "2016-03-19 00:46:56"
$epoch        = "0001-01-01 00:00:00"
$depFrom      = "e2e\$mdFile"

$srcRepoCommit = "c7470a7a2a7e5c2db7b6c3697363723741727d93"
$zhRepoCommit  = "c7d48fb50cb88276b395366cfff0ce7516591980"
$deRepoCommit  = "f961a524008aa74a7ee76b5c1153cafe004d9d06"

function SrcUrl($name) {
    return "https://github.com/OpenLocalizationTest/oltest/blob/$srcRepoCommit/e2e/$name"
}
function ZhUrl($name) {
    return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$zhRepoCommit/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$name"
}
function DeUrl($name) {
    return "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/$deRepoCommit/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$name"
}

# ---------------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Drop the existing hyperlink so it can be rebuilt against the new file name.
$ov.Range("A2").Hyperlinks.Delete()

$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"
$ov.Range("D2").Value = $overviewDate

$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"
$ov.Range("D3").Value = $overviewDate

$ov.Range("B4").Value = "Ready for handoff"
$ov.Range("C4").Value = "Ready for handoff"
$ov.Range("D4").Value = $overviewDate

$ov.Hyperlinks.Add($ov.Range("A2"), (SrcUrl $mdFile), [Type]::Missing, [Type]::Missing, $mdFile)
$ov.Hyperlinks.Add($ov.Range("A3"), (SrcUrl $png1),   [Type]::Missing, [Type]::Missing, $png1)
$ov.Hyperlinks.Add($ov.Range("A4"), (SrcUrl $png2),   [Type]::Missing, [Type]::Missing, $png2)

# ---------------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Hyperlinks.Delete()
$zh.Range("B2").Hyperlinks.Delete()
$zh.Range("D2").Hyperlinks.Delete()

$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("E2").Value = $zhDate
$zh.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H2").Value = $epoch
$zh.Range("I2").Value = "Include"

$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("E3").Value = $zhDate
$zh.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H3").Value = $epoch
$zh.Range("I3").Value = "IsDependency"
$zh.Range("J3").Value = $depFrom

$zh.Range("C4").Value = "Ready for handoff"
$zh.Range("E4").Value = $zhDate
$zh.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$zh.Range("H4").Value = $epoch
$zh.Range("I4").Value = "IsDependency"
$zh.Range("J4").Value = $depFrom

$zh.Hyperlinks.Add($zh.Range("A2"), (SrcUrl $mdFile), [Type]::Missing, [Type]::Missing, $mdFile)
$zh.Hyperlinks.Add($zh.Range("B2"), (SrcUrl $mdFile), [Type]::Missing, [Type]::Missing, ".md")
$zh.Hyperlinks.Add($zh.Range("D2"), (ZhUrl  $zhXlf),  [Type]::Missing, [Type]::Missing, $zhXlf)

$zh.Hyperlinks.Add($zh.Range("A3"), (SrcUrl $png1),   [Type]::Missing, [Type]::Missing, $png1)
$zh.Hyperlinks.Add($zh.Range("B3"), (SrcUrl $png1),   [Type]::Missing, [Type]::Missing, ".png")
$zh.Hyperlinks.Add($zh.Range("D3"), (ZhUrl  $png1Hb), [Type]::Missing, [Type]::Missing, $png1Hb)

$zh.Hyperlinks.Add($zh.Range("A4"), (SrcUrl $png2),   [Type]::Missing, [Type]::Missing, $png2)
$zh.Hyperlinks.Add($zh.Range("B4"), (SrcUrl $png2),   [Type]::Missing, [Type]::Missing, ".png")
$zh.Hyperlinks.Add($zh.Range("D4"), (ZhUrl  $png2Hb), [Type]::Missing, [Type]::Missing, $png2Hb)

# ---------------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Hyperlinks.Delete()
$de.Range("B2").Hyperlinks.Delete()
$de.Range("D2").Hyperlinks.Delete()

$de.Range("C2").Value = "Ready for handoff"
$de.Range("E2").Value = $deDate
$de.Range("E2").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H2").Value = $epoch
$de.Range("I2").Value = "Include"

$de.Range("C3").Value = "Ready for handoff"
$de.Range("E3").Value = $deDate
$de.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H3").Value = $epoch
$de.Range("I3").Value = "IsDependency"
$de.Range("J3").Value = $depFrom

$de.Range("C4").Value = "Ready for handoff"
$de.Range("E4").Value = $deDate
$de.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$de.Range("H4").Value = $epoch
$de.Range("I4").Value = "IsDependency"
$de.Range("J4").Value = $depFrom

$de.Hyperlinks.Add($de.Range("A2"), (SrcUrl $mdFile), [Type]::Missing, [Type]::Missing, $mdFile)
$de.Hyperlinks.Add($de.Range("B2"), (SrcUrl $mdFile), [Type]::Missing, [Type]::Missing, ".md")
$de.Hyperlinks.Add($de.Range("D2"), (DeUrl  $deXlf),  [Type]::Missing, [Type]::Missing, $deXlf)

$de.Hyperlinks.Add($de.Range("A3"), (SrcUrl $png1),   [Type]::Missing, [Type]::Missing, $png1)
$de.Hyperlinks.Add($de.Range("B3"), (SrcUrl $png1),   [Type]::Missing, [Type]::Missing, ".png")
$de.Hyperlinks.Add($de.Range("D3"), (DeUrl  $png1Hb), [Type]::Missing, [Type]::Missing, $png1Hb)

$de.Hyperlinks.Add($de.Range("A4"), (SrcUrl $png2),   [Type]::Missing, [Type]::Missing, $png2)
$de.Hyperlinks.Add($de.Range("B4"), (SrcUrl $png2),   [Type]::Missing, [Type]::Missing, ".png")
$de.Hyperlinks.Add($de.Range("D4"), (DeUrl  $png2Hb), [Type]::Missing, [Type]::Missing, $png2Hb)

Write-Host "Localization status report regenerated for handoff."
